# Polishing site & T-CI
# Add a website link (styled left-aligned) to the existing jcamp row,
# and append a new test entry row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45 (existing "jcamp" row): add website in column B, left-aligned ---
$ws.Range("B45").Value = "https://bryanhanson.github.io/bad_page.html"
$ws.Range("B45").HorizontalAlignment = -4131   # xlLeft

# --- Row 46 (new "test" row) ---
$ws.Range("A46").Value = "test"
$ws.Range("B46").Value = "https://bryanhanson.github.io/bad_page.html"
$ws.Range("B46").HorizontalAlignment = -4131   # xlLeft
$ws.Range("E46").Value = "Python"
$ws.Range("F46").Value = "test entry"
$ws.Range("G46").Value = "test entry"

# --- Update view: scroll/selection position ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("G46").Select()
